$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

$ws.Range("D2").Value = 0.003961923997849226
$ws.Range("E2").Value = 0.06242087110877037
$ws.Range("G2").Value = 0.004758586175739765
$ws.Range("H2").Value = 0.03974230634048581
$ws.Range("I2").Value = 0.001530169509351254
$ws.Range("J2").Value = 0.009129321202635765
$ws.Range("K2").Value = 0.00285257725045085
$ws.Range("D3").Value = 0.002390699461102486
$ws.Range("E3").Value = 0.3609774098731577
$ws.Range("G3").Value = 0.03056043619289994
$ws.Range("H3").Value = 0.22402598336339
$ws.Range("I3").Value = 0.01802518498152494
$ws.Range("J3").Value = 0.04083840968087316
$ws.Range("K3").Value = 0.01866521127521992
$ws.Range("D4").Value = 0.003317888826131821
$ws.Range("E4").Value = 0.05858613131567836
$ws.Range("G4").Value = 0.004581920802593231
$ws.Range("H4").Value = 0.03666938655078411
$ws.Range("I4").Value = 0.001605357509106398
$ws.Range("J4").Value = 0.009125441778451204
$ws.Range("K4").Value = 0.002569119445979595
$ws.Range("D5").Value = 0.00149089377373457
$ws.Range("E5").Value = 0.3604000369086862
$ws.Range("G5").Value = 0.03053943580016494
$ws.Range("H5").Value = 0.2217852910980582
$ws.Range("I5").Value = 0.01878394279628992
$ws.Range("J5").Value = 0.04089313978329301
$ws.Range("K5").Value = 0.01894107880070806
$ws.Range("C6").Value = 0
$ws.Range("E6").Value = 1.106505824718624
$ws.Range("D7").Value = 0.002617233898490667
$ws.Range("E7").Value = 0.03924175398424268
$ws.Range("G7").Value = 0.002862686756998301
$ws.Range("H7").Value = 0.02421921771019697
$ws.Range("I7").Value = 0.001118030864745378
$ws.Range("J7").Value = 0.006645085755735636
$ws.Range("K7").Value = 0.001754528842866421
$ws.Range("D8").Value = 0.002056588884443045
$ws.Range("E8").Value = 0.305586124304682
$ws.Range("G8").Value = 0.02601653430610895
$ws.Range("H8").Value = 0.1881138212047517
$ws.Range("I8").Value = 0.01473014336079359
$ws.Range("J8").Value = 0.03534559765830636
$ws.Range("K8").Value = 0.01616413472220302
$ws.Range("D9").Value = 0.002436739392578602
$ws.Range("E9").Value = 0.03652191301807761
$ws.Range("G9").Value = 0.00274158688262105
$ws.Range("H9").Value = 0.02291522035375237
$ws.Range("I9").Value = 0.001087508164346218
$ws.Range("J9").Value = 0.005839633289724588
$ws.Range("K9").Value = 0.001547982916235924
$ws.Range("D10").Value = 0.001243194565176964
$ws.Range("E10").Value = 0.2937803077511489
$ws.Range("G10").Value = 0.02501813694834709
$ws.Range("H10").Value = 0.1808811011724174
$ws.Range("I10").Value = 0.0151786208152771
$ws.Range("J10").Value = 0.03359454357996583
$ws.Range("K10").Value = 0.01547997817397118
$ws.Range("C11").Value = 0
$ws.Range("E11").Value = 0.8979644482024014
$ws.Range("D12").Value = 0.002382888458669186
$ws.Range("E12").Value = 0.04156492277979851
$ws.Range("G12").Value = 0.00315969018265605
$ws.Range("H12").Value = 0.0261453534476459
$ws.Range("I12").Value = 0.001290869433432817
$ws.Range("J12").Value = 0.006088669877499342
$ws.Range("K12").Value = 0.001958401873707771
$ws.Range("D13").Value = 0.001847608014941216
$ws.Range("E13").Value = 0.2679000240750611
$ws.Range("G13").Value = 0.0221830909140408
$ws.Range("H13").Value = 0.1667156340554357
$ws.Range("I13").Value = 0.01391408080235124
$ws.Range("J13").Value = 0.02962920162826777
$ws.Range("K13").Value = 0.01398213813081384
$ws.Range("D14").Value = 0.002360836137086153
$ws.Range("E14").Value = 0.04394010920077562
$ws.Range("G14").Value = 0.003344262484461069
$ws.Range("H14").Value = 0.0275275120511651
$ws.Range("I14").Value = 0.001386471092700958
$ws.Range("J14").Value = 0.00632206117734313
$ws.Range("K14").Value = 0.002075059805065393
$ws.Range("D15").Value = 0.001198198180645704
$ws.Range("E15").Value = 0.2847359217703342
$ws.Range("G15").Value = 0.02414836501702666
$ws.Range("H15").Value = 0.1748937289230525
$ws.Range("I15").Value = 0.01591190975159407
$ws.Range("J15").Value = 0.03171652369201183
$ws.Range("K15").Value = 0.01497953571379185
$ws.Range("C16").Value = 0
$ws.Range("E16").Value = 1.059409182518721
$ws.Range("D17").Value = 0.002766044810414314
$ws.Range("E17").Value = 0.05008273618295789
$ws.Range("G17").Value = 0.003728274255990982
$ws.Range("H17").Value = 0.03120582643896341
$ws.Range("I17").Value = 0.001357139553874731
$ws.Range("J17").Value = 0.007799938321113586
$ws.Range("K17").Value = 0.00234952662140131
$ws.Range("D18").Value = 0.002069283742457628
$ws.Range("E18").Value = 0.3067218810319901
$ws.Range("G18").Value = 0.02578407153487206
$ws.Range("H18").Value = 0.1884060804732144
$ws.Range("I18").Value = 0.01634280104190111
$ws.Range("J18").Value = 0.03542019426822662
$ws.Range("K18").Value = 0.01606227504089475
$ws.Range("D19").Value = 0.00293141882866621
$ws.Range("E19").Value = 0.05188515409827232
$ws.Range("G19").Value = 0.00397630175575614
$ws.Range("H19").Value = 0.03239032998681068
$ws.Range("I19").Value = 0.001482035033404827
$ws.Range("J19").Value = 0.008001338224858046
$ws.Range("K19").Value = 0.002349940128624439
$ws.Range("D20").Value = 0.001264994964003563
$ws.Range("E20").Value = 0.300749619025737
$ws.Range("G20").Value = 0.0252198320813477
$ws.Range("H20").Value = 0.1854279190301895
$ws.Range("I20").Value = 0.01708532124757767
$ws.Range("J20").Value = 0.03404026012867689
$ws.Range("K20").Value = 0.01527096331119537
$ws.Range("C21").Value = 0
$ws.Range("E21").Value = 1.00795480562374
$ws.Range("D22").Value = 0.002676998730748892
$ws.Range("E22").Value = 0.04325457895174623
$ws.Range("G22").Value = 0.003306585364043713
$ws.Range("H22").Value = 0.02701234305277467
$ws.Range("I22").Value = 0.00129800708964467
$ws.Range("J22").Value = 0.006706702057272196
$ws.Range("K22").Value = 0.001983508002012968
$ws.Range("D23").Value = 0.001876875292509794
$ws.Range("E23").Value = 0.2856964897364378
$ws.Range("G23").Value = 0.02377653401345015
$ws.Range("H23").Value = 0.1768475039862096
$ws.Range("I23").Value = 0.01542425900697708
$ws.Range("J23").Value = 0.0322994259186089
$ws.Range("K23").Value = 0.01486537978053093
$ws.Range("D24").Value = 0.002617469057440758
$ws.Range("E24").Value = 0.04459123313426971
$ws.Range("G24").Value = 0.003365387208759785
$ws.Range("H24").Value = 0.0278595769777894
$ws.Range("I24").Value = 0.001256630290299654
$ws.Range("J24").Value = 0.006917370017617941
$ws.Range("K24").Value = 0.002033184748142958
$ws.Range("D25").Value = 0.00120099913328886
$ws.Range("E25").Value = 0.2911525252275169
$ws.Range("G25").Value = 0.02450546575710177
$ws.Range("H25").Value = 0.178148933686316
$ws.Range("I25").Value = 0.01639682101085782
$ws.Range("J25").Value = 0.03340143710374832
$ws.Range("K25").Value = 0.0150182475335896
$ws.Range("C26").Value = 0
$ws.Range("E26").Value = 0.8913419921882451
$ws.Range("D27").Value = 0.002881017979234457
$ws.Range("E27").Value = 0.04731297260150313
$ws.Range("G27").Value = 0.003563164547085762
$ws.Range("H27").Value = 0.02966500939801335
$ws.Range("I27").Value = 0.00131884329020977
$ws.Range("J27").Value = 0.007273943442851305
$ws.Range("K27").Value = 0.002179708518087864
$ws.Range("D28").Value = 0.002048211079090834
$ws.Range("E28").Value = 0.3053763858042657
$ws.Range("G28").Value = 0.02566413339227438
$ws.Range("H28").Value = 0.1888218046166003
$ws.Range("I28").Value = 0.01568729383870959
$ws.Range("J28").Value = 0.03470656583085656
$ws.Range("K28").Value = 0.01594782778993249
$ws.Range("D29").Value = 0.002732870448380709
$ws.Range("E29").Value = 0.04710490815341473
$ws.Range("G29").Value = 0.003601891826838255
$ws.Range("H29").Value = 0.02947240518406034
$ws.Range("I29").Value = 0.001363600417971611
$ws.Range("J29").Value = 0.007241168897598982
$ws.Range("K29").Value = 0.002115057408809662
$ws.Range("D30").Value = 0.001279656123369932
$ws.Range("E30").Value = 0.3061636821366847
$ws.Range("G30").Value = 0.02588624712079763
$ws.Range("H30").Value = 0.1882273947820067
$ws.Range("I30").Value = 0.01667132312431931
$ws.Range("J30").Value = 0.03472918085753918
$ws.Range("K30").Value = 0.01593796070665121
$ws.Range("C31").Value = 0
$ws.Range("E31").Value = 0.9926352506503463
